$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Row_Number" column at G, pushing the old
# TextBlob_Sentiment_Label/Score columns from G/H to H/I.
$ws.Columns("G").Insert(-4161)
$ws.Range("G1").Value = "Row_Number"

# Fill the new Row_Number column for the 21 remaining data rows (2-22)
# with the sequence 1..21.
for ($i = 2; $i -le 22; $i++) {
    $ws.Cells.Item($i, 7).Value = $i - 1
}

# Drop the last four rows (the extra Adidas-related news rows that were
# removed), shifting everything below them up.
$ws.Range("A23:I26").Delete(-4162)
